# Add a new row for customer 79174408 with 0 points (matching the
# existing "phone, birthday, total_points" table layout).
#
# Row 11: A11 = "79174408" (text, like the other phone numbers end up
#         being stored), B11 = "" (blank birthday, same pattern as the
#         other rows with no birthday on file), C11 = 0 (total_points).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the numeric-looking phone number to be
# stored as text instead of being auto-converted to a number, matching
# column A's existing text-like phone values.
$ws.Cells.Item(11, 1).Value = "'79174408"
$ws.Cells.Item(11, 1).Style = "Normal"

# Blank birthday cell - written the same way (quote-prefixed empty
# text) as the other "no birthday on file" rows so it ends up as an
# empty text value rather than a totally missing cell.
$ws.Cells.Item(11, 2).Value = "'"
$ws.Cells.Item(11, 2).Style = "Normal"

# total_points starts at 0 for the new customer.
$ws.Cells.Item(11, 3).Value = 0
